$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new data points (data type specific behaviors test rows)
$ws.Range("L7").Value = 1
$ws.Range("B13").Value = 1

# Move the active selection the way Excel would leave it after this entry
$ws.Range("B14").Select()
